$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07107759697255688
$ws.Range("C2").Value = 0.6196869661694117
$ws.Range("B3").Value = 0.07469256291730456
$ws.Range("C3").Value = 0.859566932473949
